# Update the "Cookie Types" workbook for the Power Pivot exercise:
#  - rescale the Revenue/Cost per cookie figures from (assumed) dollars to
#    Indonesian Rupiah-sized integers (x1000)
#  - format those columns as Indonesian Rupiah currency ("Rp")
#  - move the active selection
#  - force portrait page orientation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rescale the data values in B2:C7 (Revenue Per Cookie / Cost Per Cookie)
$ws.Range("B2").Value2 = 5000
$ws.Range("C2").Value2 = 2000

$ws.Range("B3").Value2 = 1000
$ws.Range("C3").Value2 = 500

$ws.Range("B4").Value2 = 5000
$ws.Range("C4").Value2 = 2200

$ws.Range("B5").Value2 = 4000
$ws.Range("C5").Value2 = 1500

$ws.Range("B6").Value2 = 3000
$ws.Range("C6").Value2 = 1250

$ws.Range("B7").Value2 = 6000
$ws.Range("C7").Value2 = 2750

# 2. Apply the Indonesian Rupiah currency number format to the data columns
$ws.Range("B2:C7").NumberFormat = '_-"Rp"* #,##0_-;\-"Rp"* #,##0_-;_-"Rp"* "-"_-;_-@_-'

# 3. Update the selected cell
$ws.Range("A12").Select()

# 4. Force the print page orientation to portrait
$ws.PageSetup.Orientation = 1
